{"js": "// Update the date paragraph and the 25 multiplication answers in the single\n// table (5 \"data\" rows x 5 columns; the rows in between are blank spacers).\nconst body = context.document.body;\n\n// --- 1. Title paragraph: \"2024-09-21 Saturday\" -> \"2024-09-22 Sunday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2024-09-22 Sunday\", \"Replace\");\n\n// --- 2. Table cell values ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only these five rows actually contain answers; the others are blank\n// spacer rows. Values are given left-to-right, top-to-bottom matching the\n// diff order.\nconst newValues = {\n  0: [\"959\u00d73=2877\", \"687\u00d73=2061\", \"267\u00d79=2403\", \"862\u00d75=4310\", \"163\u00d79=1467\"],\n  4: [\"829\u00d72=1658\", \"754\u00d79=6786\", \"598\u00d73=1794\", \"422\u00d74=1688\", \"432\u00d79=3888\"],\n  9: [\"310\u00d74=1240\", \"385\u00d78=3080\", \"142\u00d72=284\", \"687\u00d74=2748\", \"111\u00d79=999\"],\n  14: [\"196\u00d76=1176\", \"611\u00d73=1833\", \"131\u00d74=524\", \"597\u00d77=4179\", \"251\u00d79=2259\"],\n  19: [\"991\u00d76=5946\", \"601\u00d73=1803\", \"313\u00d76=1878\", \"234\u00d73=702\", \"868\u00d79=7812\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = Number(rowIndex);\n  const rowVals = newValues[rowIndex];\n  for (let c = 0; c < rowVals.length; c++) {\n    table.getCell(r, c).value = rowVals[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 multiplication answers in the single\n# table (5 \"data\" rows x 5 columns; the rows in between are blank spacers).\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"2024-09-21 Saturday\" -> \"2024-09-22 Sunday\" ---\n$d.Paragraphs(1).Range.Text = \"2024-09-22 Sunday\"\n\n# --- 2. Table cell values ---\n# Only these five (1-based) rows actually contain answers; the others are\n# blank spacer rows. Values are given left-to-right matching the diff order.\n$t = $d.Tables(1)\n\n$newValues = @{\n    1  = @(\"959\u00d73=2877\", \"687\u00d73=2061\", \"267\u00d79=2403\", \"862\u00d75=4310\", \"163\u00d79=1467\")\n    5  = @(\"829\u00d72=1658\", \"754\u00d79=6786\", \"598\u00d73=1794\", \"422\u00d74=1688\", \"432\u00d79=3888\")\n    10 = @(\"310\u00d74=1240\", \"385\u00d78=3080\", \"142\u00d72=284\",  \"687\u00d74=2748\", \"111\u00d79=999\")\n    15 = @(\"196\u00d76=1176\", \"611\u00d73=1833\", \"131\u00d74=524\",  \"597\u00d77=4179\", \"251\u00d79=2259\")\n    20 = @(\"991\u00d76=5946\", \"601\u00d73=1803\", \"313\u00d76=1878\", \"234\u00d73=702\",  \"868\u00d79=7812\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowVals = $newValues[$rowIndex]\n    for ($c = 0; $c -lt $rowVals.Count; $c++) {\n        $t.Cell($rowIndex, $c + 1).Range.Text = $rowVals[$c]\n    }\n}\n"}
